$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 183 ("Segunda" quality, week-on-week entry),
# pushing the existing rows 183-266 down to 184-267 and extending the
# sheet's used range from A1:R266 to A1:R267.
$ws.Rows.Item(183).Insert()

$ws.Range("A183").Value = 5
$ws.Range("B183").Value = "Macroferia Regional de Talca"
$ws.Range("C183").Value = "Maule"
$ws.Range("D183").Value = 44609
$ws.Range("D183").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 100114014
$ws.Range("G183").Value = "Betarraga"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Segunda"
$ws.Range("J183").Value = 3000
$ws.Range("K183").Value = 700
$ws.Range("L183").Value = 700
$ws.Range("M183").Value = 700
$ws.Range("N183").Value = "`$/paquete 5 unidades"
$ws.Range("O183").Value = "Región del Maule"
$ws.Range("P183").Value = 140
$ws.Range("Q183").Value = 5
$ws.Range("R183").Value = "Hortaliza"
